$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format from an existing date cell (A18) down into the new rows
$ws.Range("A18").Copy()
$ws.Range("A19:A22").PasteSpecial(-4122)  # xlPasteFormats

# Row 19
$ws.Range("A19").Value = 43399
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = "Meeting"

# Row 20
$ws.Range("A20").Value = 43404
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = "Fietssimulatie bijwerken"

# Row 21
$ws.Range("A21").Value = 43406
$ws.Range("B21").Value = 2
$ws.Range("C21").Formula = "=C20"

# Row 22
$ws.Range("A22").Value = 43407
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = "Preprocessing"

# Update selection to match diff
$ws.Range("B24").Select()

$wb.Save()
